$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (certificate holder) is replaced with the real registrant's details
$ws.Range("A2").Value = "RAGHUVARAN E"
$ws.Range("B2").Value = "ASSISTANT PROFESSOR"
$ws.Range("F2").Value = "SRI VENKATESWARA COLLEGE OF ENGINEERING"

# Column C header: "Date" -> "Mail"
$ws.Range("C1").Value = "Mail"
$ws.Range("C2").Value = "eraghu21@gmail.com"

# Turn the e-mail address into a real mailto: hyperlink (Excel auto-applies
# the built-in "Hyperlink" cell style to C2 when this runs)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:eraghu21@gmail.com")

# Re-apply bold to the header row so Excel's style table collapses the
# header cells back onto a single, shared cell-format index
$ws.Range("A1:F1").Font.Bold = $true
